$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set NumberFormat to Text for D-column cells being updated, to preserve
# literal strings such as "225.75" or "1.00" instead of Excel auto-converting
# them to numeric values.
$dCells = @("D2", "D3", "D5", "D6", "D8", "D10", "D11", "D12", "D14", "D16", "D17", "D18", "D19", "D21", "D25", "D26", "D27", "D32", "D35", "D38", "D40", "D42", "D44", "D45", "D46", "D47", "D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
$ws.Range("D2").Value = '34.501.10'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.809.74'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '225.75'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").Value = '0.598'
$ws.Range("E6").Value = '  +2.56%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '36.21'
$ws.Range("E8").Value = '  +3.79%  '
$ws.Range("E9").Value = '  -3.00%  '
$ws.Range("D10").Value = '0.0681'
$ws.Range("E10").Value = '  -2.04%  '
$ws.Range("D11").Value = '0.0967'
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("D12").Value = '2.070.68'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").Value = '1.799.56'
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D16").Value = '34.473.12'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '4.41'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '68.40'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").Value = '242.24'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("E20").Value = '  -3.05%  '
$ws.Range("D21").Value = '11.22'
$ws.Range("E21").Value = '  -2.87%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("E23").Value = '  -1.72%  '
$ws.Range("E24").Value = '  +5.13%  '
$ws.Range("D25").Value = '171.54'
$ws.Range("E25").Value = '  -1.85%  '
$ws.Range("D26").Value = '7.87'
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("D27").Value = '17.25'
$ws.Range("E27").Value = '  +2.61%  '
$ws.Range("E28").Value = '  +1.62%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("E31").Value = '  -1.40%  '
$ws.Range("D32").Value = '3.90'
$ws.Range("E32").Value = '  -2.98%  '
$ws.Range("E33").Value = '  -2.39%  '
$ws.Range("D35").Value = '1.363.26'
$ws.Range("E35").Value = '  -2.46%  '
$ws.Range("E36").Value = '  -4.30%  '
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("D38").Value = '2.36'
$ws.Range("E38").Value = '  -6.75%  '
$ws.Range("E39").Value = '  -2.27%  '
$ws.Range("D40").Value = '2.43'
$ws.Range("E40").Value = '  +0.60%  '
$ws.Range("E41").Value = '  -1.86%  '
$ws.Range("D42").Value = '80.89'
$ws.Range("E42").Value = '  -3.08%  '
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").Value = '1.16'
$ws.Range("E44").Value = '  +3.97%  '
$ws.Range("D45").Value = '13.39'
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("D46").Value = '0.0499'
$ws.Range("E46").Value = '  -2.52%  '
$ws.Range("D47").Value = '1.970.98'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("E48").Value = '  -2.77%  '
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").Value = '102.41'
$ws.Range("E50").Value = '  -2.50%  '
$ws.Range("E51").Value = '  -5.93%  '

# Restore default styling on the D-column cells (undo the temporary text format)
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
